$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "@yassine"

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "pierrick"

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "@yassine2"
